# Applies the Thu Jul 27 04:13:45 UTC 2023 cryptos-list refresh:
# updated prices/volume deltas and re-ranked rows 21-51 to match
# the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''29.423.71'
$ws.Range("E2").Value = '  +0.65%  '

# Row 3
$ws.Range("D3").Value = '''1.877.27'
$ws.Range("E3").Value = '  +1.04%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").Value = '''0.7184'
$ws.Range("E5").Value = '  +1.21%  '

# Row 6
$ws.Range("D6").Value = '''240.17'
$ws.Range("E6").Value = '  +0.65%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("D8").Value = '''0.07831'
$ws.Range("E8").Value = '  -2.12%  '

# Row 9
$ws.Range("D9").Value = '''0.3114'
$ws.Range("E9").Value = '  +2.63%  '

# Row 10
$ws.Range("D10").Value = '''24.96'
$ws.Range("E10").Value = '  +6.28%  '

# Row 11
$ws.Range("D11").Value = '''0.08248'
$ws.Range("E11").Value = '  +0.54%  '

# Row 12
$ws.Range("D12").Value = '''1.887.24'
$ws.Range("E12").Value = '  +2.64%  '

# Row 13
$ws.Range("D13").Value = '''0.7272'
$ws.Range("E13").Value = '  +3.25%  '

# Row 14
$ws.Range("D14").Value = '''5.287'
$ws.Range("E14").Value = '  +2.12%  '

# Row 15
$ws.Range("D15").Value = '''91.31'
$ws.Range("E15").Value = '  +1.89%  '

# Row 16
$ws.Range("D16").Value = '''29.455.72'
$ws.Range("E16").Value = '  +0.88%  '

# Row 17
$ws.Range("D17").Value = '''5.927'
$ws.Range("E17").Value = '  +1.81%  '

# Row 18
$ws.Range("D18").Value = '''245.33'
$ws.Range("E18").Value = '  +3.03%  '

# Row 19
$ws.Range("D19").Value = '''0.000007882'
$ws.Range("E19").Value = '  +0.11%  '

# Row 20
$ws.Range("E20").Value = '  +0.25%  '

# Row 21
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '''1.0000'
$ws.Range("E21").Value = '  +0.11%  '

# Row 22
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").Value = '''7.962'
$ws.Range("E22").Value = '  +6.62%  '

# Row 23
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").Value = '''0.9997'
$ws.Range("E23").Value = '  -0.04%  '

# Row 24
$ws.Range("B24").Value = 'Stellar'
$ws.Range("C24").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D24").Value = '''0.1583'
$ws.Range("E24").Value = '  +9.58%  '

# Row 25
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '''163.86'
$ws.Range("E25").Value = '  +0.58%  '

# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '''9.045'
$ws.Range("E26").Value = '  +1.34%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''18.31'

# Row 28
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '''1.365'
$ws.Range("E28").Value = '  -4.41%  '

# Row 29
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '''1.484'
$ws.Range("E29").Value = '  +0.02%  '

# Row 30
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = '''4.383'
$ws.Range("E30").Value = '  +0.36%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '''4.153'
$ws.Range("E31").Value = '  +3.44%  '

# Row 32
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '''0.05280'
$ws.Range("E32").Value = '  +1.46%  '

# Row 33
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").Value = '''1.943'
$ws.Range("E33").Value = '  +1.28%  '

# Row 34
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '''1.200'
$ws.Range("E34").Value = '  +3.44%  '

# Row 35
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '''0.7215'
$ws.Range("E35").Value = '  +1.53%  '

# Row 36
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '''2.679'
$ws.Range("E36").Value = '  +0.13%  '

# Row 37
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '''0.01866'
$ws.Range("E37").Value = '  +0.63%  '

# Row 38
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '''1.240.00'
$ws.Range("E38").Value = '  +9.34%  '

# Row 39
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '''2.726'
$ws.Range("E39").Value = '  +0.24%  '

# Row 40
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '''0.9048'
$ws.Range("E40").Value = '  -2.73%  '

# Row 41
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '''73.65'
$ws.Range("E41").Value = '  +5.13%  '

# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''6.099'
$ws.Range("E42").Value = '  +3.94%  '

# Row 43
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = '''1.000'
$ws.Range("E43").Value = '  +0.04%  '

# Row 44
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '''103.36'
$ws.Range("E44").Value = '  +0.80%  '

# Row 45
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '''0.5336'
$ws.Range("E45").Value = '  -0.01%  '

# Row 46
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '''0.00000000121'
$ws.Range("E46").Value = '  +0.91%  '

# Row 47
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '''1.760'
$ws.Range("E47").Value = '  -0.23%  '

# Row 48
$ws.Range("D48").Value = '''2.911'
$ws.Range("E48").Value = '  +12.12%  '

# Row 49
$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").Value = '''0.4330'
$ws.Range("E49").Value = '  +1.73%  '

# Row 50
$ws.Range("D50").Value = '''9.271'
$ws.Range("E50").Value = '  +1.06%  '

# Row 51
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = '''7.078'
$ws.Range("E51").Value = '  +1.63%  '

